$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Series"-level rows appended below the header, describing MCH121-1 and
# MCH121-2. Columns: A=identifier, E=levelOfDescription, F=extentAndMedium,
# G=notes (matching the header row already in the sheet).
$rows = @(
    @{ Row = 2; Identifier = "MCH121-1" },
    @{ Row = 3; Identifier = "MCH121-2" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Style every cell in the row (including the blank placeholder cells)
    # with the plain, non-bold Calibri 10pt body font used for data rows.
    foreach ($col in @("A","C","D","E","F","G","H")) {
        $cell = $ws.Range("$col$rowNum")
        $cell.Font.Name = "Calibri"
        $cell.Font.Size = 10
        $cell.Font.ThemeColor = 1
    }

    $ws.Range("A$rowNum").Value = $r.Identifier
    $ws.Range("E$rowNum").Value = "Series"
    $ws.Range("F$rowNum").Value = "1 Box"
    $ws.Range("G$rowNum").Value = "LOCATION: 21D | GRAP COUNT NUMER: NONE"
}

# Match the saved selection from the edit (cell C12).
$ws.Range("C12").Select() | Out-Null
